$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计")
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item(3)
$ws = $wb.Worksheets.Add($null, $q4)
$ws.Name = "2022-Q1"

# Re-use the existing header / index-column formatting (same style used by
# the other quarterly sheets) instead of inventing new style records.
$q4.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q4.Range("A2").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)   # xlPasteFormats

$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row
$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# Data rows: A(idx), 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$data = @(
    @(0, "012719", "华夏新兴经济一年持有期混合型证券投资基金A", "21.62", "91.19", "3.03", "0.6551", 10),
    @(1, "012421", "华夏优加生活混合A", "8.67", "92.98", "2.42", "0.2098", 10),
    @(2, "160921", "大成多策略混合(LOF)", "1.13", "79.19", "6.82", "0.0771", 2),
    @(3, "160519", "博时睿利事件驱动灵活配置混合", "1.40", "80.50", "3.37", "0.0472", 4),
    @(4, "012720", "华夏新兴经济一年持有期混合型证券投资基金C", "1.04", "91.19", "3.03", "0.0315", 10),
    @(5, "012422", "华夏优加生活混合C", "0.17", "92.98", "2.42", "0.0041", 10)
)

$r = 2
foreach ($row in $data) {
    # Column A keeps its pasted-in style and is a genuine number.
    $ws.Cells.Item($r,1).Value = $row[0]

    # Columns B..G must stay plain text (fund codes keep leading zeros, and
    # the numeric-looking figures are stored as text in the source data), so
    # force a text number-format before assignment, then drop back to the
    # built-in "Normal" style so no stray formatting is left behind.
    for ($col = 2; $col -le 7; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col - 1]
        $cell.Style = "Normal"
    }

    # Column H (rank) is a genuine number.
    $ws.Cells.Item($r,8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (Total) sheet: add a 2022-Q1 summary row at the top of
#    the data (row 2), pushing the existing quarters down and renumbering the
#    index column (A).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 6
$total.Cells.Item(2,4).Value = 1.02

# Give the new A2 the same styling as the other index cells in column A.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# Renumber the index column 0..3 top to bottom.
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3

# Keep the original active sheet selected.
$wb.Worksheets.Item(1).Activate()
